# Updates cryptos list figures (Price / Volume(1h) columns) and one
# coin rename (row 51: MultiversX -> BitcoinSV) per the scheduled scrape.
#
# Every cell is written through a temporary formula (`="<text>"`) that is
# then flattened to a plain value via Copy/PasteSpecial(xlPasteValues).
# Writing straight to `.Value` lets Excel auto-convert number-looking
# strings (e.g. "314.90" -> 314.9, dropping the trailing zero and the
# original text formatting) which does not match the source data, where
# every cell in these columns is stored as literal text. The formula round
# trip keeps the exact text (incl. padding spaces in the % column) without
# forcing a NumberFormat change on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

Set-TextValue $ws.Cells.Item(2, 4) "41.683.05"
Set-TextValue $ws.Cells.Item(2, 5) "  +0.31%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.465.05"
Set-TextValue $ws.Cells.Item(3, 5) "  -1.16%  "
Set-TextValue $ws.Cells.Item(4, 5) "  +0.64%  "
Set-TextValue $ws.Cells.Item(5, 4) "314.90"
Set-TextValue $ws.Cells.Item(5, 5) "  +0.54%  "
Set-TextValue $ws.Cells.Item(6, 4) "92.27"
Set-TextValue $ws.Cells.Item(6, 5) "  -2.10%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.549"
Set-TextValue $ws.Cells.Item(7, 5) "  +0.31%  "
Set-TextValue $ws.Cells.Item(8, 5) "  +0.56%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.512"
Set-TextValue $ws.Cells.Item(9, 5) "  +2.57%  "
Set-TextValue $ws.Cells.Item(10, 4) "32.44"
Set-TextValue $ws.Cells.Item(10, 5) "  -1.77%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.0842"
Set-TextValue $ws.Cells.Item(11, 5) "  +7.11%  "
Set-TextValue $ws.Cells.Item(12, 5) "  +0.57%  "
Set-TextValue $ws.Cells.Item(13, 4) "2.842.85"
Set-TextValue $ws.Cells.Item(13, 5) "  -1.32%  "
Set-TextValue $ws.Cells.Item(14, 4) "6.87"
Set-TextValue $ws.Cells.Item(14, 5) "  -0.36%  "
Set-TextValue $ws.Cells.Item(15, 4) "15.79"
Set-TextValue $ws.Cells.Item(15, 5) "  +1.66%  "
Set-TextValue $ws.Cells.Item(16, 4) "2.480.84"
Set-TextValue $ws.Cells.Item(16, 5) "  -0.88%  "
Set-TextValue $ws.Cells.Item(17, 4) "0.778"
Set-TextValue $ws.Cells.Item(17, 5) "  +2.30%  "
Set-TextValue $ws.Cells.Item(18, 4) "41.647.26"
Set-TextValue $ws.Cells.Item(18, 5) "  -0.15%  "
Set-TextValue $ws.Cells.Item(19, 5) "  +2.25%  "
Set-TextValue $ws.Cells.Item(20, 5) "  +2.54%  "
Set-TextValue $ws.Cells.Item(21, 4) "70.67"
Set-TextValue $ws.Cells.Item(21, 5) "  +0.16%  "
Set-TextValue $ws.Cells.Item(22, 4) "11.41"
Set-TextValue $ws.Cells.Item(22, 5) "  +1.34%  "
Set-TextValue $ws.Cells.Item(23, 4) "238.83"
Set-TextValue $ws.Cells.Item(23, 5) "  +0.90%  "
Set-TextValue $ws.Cells.Item(24, 4) "2.71"
Set-TextValue $ws.Cells.Item(24, 5) "  -0.39%  "
Set-TextValue $ws.Cells.Item(25, 5) "  +0.40%  "
Set-TextValue $ws.Cells.Item(26, 5) "  +0.07%  "
Set-TextValue $ws.Cells.Item(27, 4) "24.45"
Set-TextValue $ws.Cells.Item(27, 5) "  -1.02%  "
Set-TextValue $ws.Cells.Item(28, 5) "  +0.39%  "
Set-TextValue $ws.Cells.Item(29, 5) "  +0.41%  "
Set-TextValue $ws.Cells.Item(30, 4) "35.17"
Set-TextValue $ws.Cells.Item(30, 5) "  -3.44%  "
Set-TextValue $ws.Cells.Item(31, 4) "155.68"
Set-TextValue $ws.Cells.Item(31, 5) "  +0.72%  "
Set-TextValue $ws.Cells.Item(32, 4) "5.49"
Set-TextValue $ws.Cells.Item(32, 5) "  +1.00%  "
Set-TextValue $ws.Cells.Item(33, 4) "2.58"
Set-TextValue $ws.Cells.Item(33, 5) "  +0.16%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.0761"
Set-TextValue $ws.Cells.Item(34, 5) "  +0.36%  "
Set-TextValue $ws.Cells.Item(35, 5) "  -0.59%  "
Set-TextValue $ws.Cells.Item(36, 4) "17.44"
Set-TextValue $ws.Cells.Item(36, 5) "  -5.17%  "
Set-TextValue $ws.Cells.Item(37, 5) "  -2.52%  "
Set-TextValue $ws.Cells.Item(38, 5) "  +0.69%  "
Set-TextValue $ws.Cells.Item(39, 5) "  +0.87%  "
Set-TextValue $ws.Cells.Item(40, 5) "  -2.59%  "
Set-TextValue $ws.Cells.Item(41, 4) "3.94"
Set-TextValue $ws.Cells.Item(41, 5) "  -5.45%  "
Set-TextValue $ws.Cells.Item(42, 5) "  +0.55%  "
Set-TextValue $ws.Cells.Item(43, 4) "1.974.77"
Set-TextValue $ws.Cells.Item(43, 5) "  +1.03%  "
Set-TextValue $ws.Cells.Item(44, 5) "  -1.23%  "
Set-TextValue $ws.Cells.Item(45, 4) "18.79"
Set-TextValue $ws.Cells.Item(45, 5) "  -4.66%  "
Set-TextValue $ws.Cells.Item(46, 4) "2.94"
Set-TextValue $ws.Cells.Item(46, 5) "  -1.63%  "
Set-TextValue $ws.Cells.Item(47, 4) "9.00"
Set-TextValue $ws.Cells.Item(47, 5) "  +1.85%  "
Set-TextValue $ws.Cells.Item(48, 4) "2.699.85"
Set-TextValue $ws.Cells.Item(48, 5) "  -1.33%  "
Set-TextValue $ws.Cells.Item(49, 4) "96.81"
Set-TextValue $ws.Cells.Item(49, 5) "  +0.00%  "
Set-TextValue $ws.Cells.Item(50, 4) "66.93"
Set-TextValue $ws.Cells.Item(50, 5) "  -1.19%  "
Set-TextValue $ws.Cells.Item(51, 2) "BitcoinSV"
Set-TextValue $ws.Cells.Item(51, 3) "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Cells.Item(51, 4) "72.26"
Set-TextValue $ws.Cells.Item(51, 5) "  -1.72%  "

$excel.CutCopyMode = 0

